# Update the symptom_group table: insert "Administrative" row, drop
# "Symptom – Musculoskeletal" row, and refresh all percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Administrative",             2,    1.9, 1.1),
    @("Diseases (patient-stated)",  6.4,  2.4, 3.6),
    @("Injuries & adverse effects", 18.4, 7.5, 8.5),
    @("Other",                      6.4,  4.4, 7.4),
    @("Symptom – Circulatory",      8,    5.5, 8.199999999999999),
    @("Symptom – Digestive",        10.9, 4.8, 10),
    @("Symptom – General",          3.6,  2.4, 6.2),
    @("Symptom – Nervous",          11.7, 12.3, 10),
    @("Symptom – Respiratory",      27.4, 57.8, 40.5),
    @("Symptom – Skin/Hair/Nails",  2.1,  0,   3),
    @("Uncodable/Unknown",          3.1,  1.1, 1.4)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}
